$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E36").Value = 1
$ws.Range("E38").Value = 1
$ws.Range("E39").Value = 1
$ws.Range("E40").Value = 1

$ws.Range("E41").Select()
